$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price & volume figures refreshed; some rows
# shifted because PaxosStandard dropped out of the list and Elrond was added).

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '24.980.01'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -3.77%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.636.11'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -5.80%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '231.40'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -6.26%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4703'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -6.77%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2528'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -7.47%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06048'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -2.30%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.06993'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.72%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.635.42'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -5.94%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '14.20'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -7.06%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.280'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -9.93%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.5603'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -14.78%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '72.95'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -6.04%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '24.974.29'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -3.85%  '
$ws.Range('E19').Value = '  -6.22%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.000006523'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -4.84%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.848.42'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -5.81%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.241'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -8.19%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '8.415'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -4.50%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.163'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -5.33%  '
$ws.Range('E25').Value = '  -1.48%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '14.77'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.32%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '1.365'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -6.53%  '
$ws.Range('E28').Value = '  -2.34%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.621'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -9.55%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '3.852'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -3.44%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.07525'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -7.55%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.480'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -6.33%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.9993'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.04%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.04199'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -11.79%  '
$ws.Range('E35').Value = '  -3.40%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.9246'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -7.31%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.5837'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -4.52%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.563'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -6.47%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.8680'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.37%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.01468'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -9.03%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '97.30'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.52%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.745'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -10.62%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.3625'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -7.97%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '4.609'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -8.21%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.05192'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.1083'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -8.51%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '6.003'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -6.06%  '
$ws.Range('B49').Value = 'TrueUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '28.45'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -8.09%  '
